$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the actual-result outcome ("Didn't login, Test passed") for every
# data row in column C, alongside the existing Email / Pass / Expected
# Result columns (A / B / D).
$ws.Range("C2:C7").Value = "Didn't login, Test passed"
